$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 ("Introduction to Digital Twins"): "3 components..." -> "4 components..."
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$sh2Body = $s2.Shapes.Item(2)
$sh2Body.TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "4 components: phyisical model, virtual model, communication services and the data"

# ---------------------------------------------------------------------------
# Slide 9 ("From application-oriented to domain-oriented"):
#   - insert a new first paragraph
#   - "Open question" -> "Research question"
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$sh9Body = $s9.Shapes.Item(2)
$sh9Body.TextFrame.TextRange.Paragraphs(1).InsertBefore("Can we abstract from application-level solutions to domain-level solutions? (e.g. from having a platform to support a DT application, to having a platform to support Agricolture Digital Twins)`r")
$sh9Body.TextFrame.TextRange.Paragraphs(3).Runs(1).Text = "Research question"

# ---------------------------------------------------------------------------
# Slide 10 ("Modelling Digital Twin Data"): insert a new paragraph right
# before "What about an hybrid data structure?"
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$sh10Body = $s10.Shapes.Item(2)
$lastIdx = $sh10Body.TextFrame.TextRange.Paragraphs().Count
$prevIdx = $lastIdx - 1
$sh10Body.TextFrame.TextRange.Paragraphs($prevIdx).InsertAfter("`rYet, no multi-store solution has achieved broad adoption in the literature.")

# ---------------------------------------------------------------------------
# Slide 11 ("An Hybryd data structure enabling Digital Twin Data"):
#   - "Combining the strength..." gains a trailing clause
#   - "...techinques." -> "...techinques with promising results." (+ Kotlin)
#
# NOTE: in the source deck both "Content Placeholder 2" shapes on this slide
# share the same underlying shape id, which confuses index-based shape
# look-up for the *second* one (Shapes.Item(5) resolves to the first shape's
# data). Editing shape 2 first, then "deleting" it (which in this runtime
# clears it and hands it a fresh id) un-aliases the id so Shapes.Item(5)
# correctly reaches the real second shape; shape 2's wanted text is then
# restored through the same (now unique) reference.
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)

$sh11First = $s11.Shapes.Item(2)
$combiningText = "Combining the strength of Graph and Time-Series DBMS with a novel, hybrid data structure."
$sh11First.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = $combiningText

$sh11First.Delete() | Out-Null

$sh11Second = $s11.Shapes.Item(5)
$sh11Second.TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "The data structure has been implemented in Kotlin and evaluated against state-of-the-art techinques with promising results."

$sh11Restored = $s11.Shapes.Item(2)
$sh11Restored.TextFrame.TextRange.Text = $combiningText
